# Swap the "Name" and "Email" columns (A and B) for all data rows,
# then update the active selection to D8 to match the saved cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 1; $row -le 5; $row++) {
    $a = $ws.Cells.Item($row, 1).Value2
    $b = $ws.Cells.Item($row, 2).Value2
    $ws.Cells.Item($row, 1).Value2 = $b
    $ws.Cells.Item($row, 2).Value2 = $a
}

$ws.Range("D8").Select()
